# Applies:
#  1. Sheet "BME_BCCW" (sheet1): extend the three conditionalFormatting
#     sqref ranges from row 10000 to row 1000000 (no new rules).
#  2. Sheet "BME_DI_BCCW" (sheet2): extend the three conditionalFormatting
#     sqref ranges the same way, add one duplicate cfRule per block, add
#     a new data row (row 17) mirroring row 16 but with Qty (F) = 2, and
#     bump the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "BME_BCCW" -- conditional formatting range bump only.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$rngAll1 = $ws1.Range("A1:Z1000000")
$fcs1 = $rngAll1.FormatConditions
$count1 = $fcs1.Count
$snapshot1 = @()
for ($i = 1; $i -le $count1; $i++) {
    $snapshot1 += $fcs1.Item($i)
}
foreach ($fc in $snapshot1) {
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq "`$G`$10:`$G`$10000") {
        $fc.ModifyAppliesToRange($ws1.Range("G10:G1000000"))
    } elseif ($addr -eq "`$H`$10:`$K`$10000") {
        $fc.ModifyAppliesToRange($ws1.Range("H10:K1000000"))
    } elseif ($addr -eq "`$K`$9:`$K`$10000") {
        $fc.ModifyAppliesToRange($ws1.Range("K9:K1000000"))
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "BME_DI_BCCW"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# -- add new row 17, mirroring row 16's content/formulas but Qty (F) = 2
#
# Row 16 carries direct "s=4" formatting (a border + general number
# format). Priming row 17 with that formatting *before* typing the
# formulas avoids Excel's "auto-apply a number format on formula entry"
# heuristic minting a brand-new, unused cellXfs entry (it instead reuses
# the existing currency+border style already on row 16) -- then we
# re-paste the formats once more so every cell in the row (including the
# formula cells, which the heuristic nudges onto that currency style)
# ends up back on the same "s=4" style as row 16.
$ws2.Range("A16:L16").Copy()
$ws2.Range("A17:L17").PasteSpecial(-4122)

$ws2.Range("A17").Value = "PHSA"
$ws2.Range("B17").Value = "IMAG"
$ws2.Range("C17").Value = "BCCH"
$ws2.Range("D17").Value = 0
$ws2.Range("E17").Value = "ANESTHESIA UNIT VAPORIZERS"
$ws2.Range("F17").Value = 2
$ws2.Range("G17").Value = 4.949354005167959
$ws2.Range("H17").Formula = "=B10*G17"
$ws2.Range("I17").Formula = "=B11*G17"
$ws2.Range("J17").Value = 0
$ws2.Range("K17").Formula = "=SUM(H17, I17, J17)"
$ws2.Range("L17").Formula = "=K17*F17"

$ws2.Range("A16:L16").Copy()
$ws2.Range("A17:L17").PasteSpecial(-4122)

# -- extend the three conditionalFormatting sqref ranges + duplicate the
#    rule in each block (same formula/type as its siblings, new priority)
$rngAll2 = $ws2.Range("A1:Z1000000")
$fcs2 = $rngAll2.FormatConditions
$count2 = $fcs2.Count
$snapshot2 = @()
for ($i = 1; $i -le $count2; $i++) {
    $snapshot2 += $fcs2.Item($i)
}
foreach ($fc in $snapshot2) {
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq "`$G`$10:`$G`$10000") {
        $fc.ModifyAppliesToRange($ws2.Range("G10:G1000000"))
    } elseif ($addr -eq "`$H`$10:`$L`$10000") {
        $fc.ModifyAppliesToRange($ws2.Range("H10:L1000000"))
    } elseif ($addr -eq "`$L`$9:`$L`$10000") {
        $fc.ModifyAppliesToRange($ws2.Range("L9:L1000000"))
    }
}

$newG = $ws2.Range("G10:G1000000").FormatConditions.Add(2, [System.Reflection.Missing]::Value, "=LEN(TRIM(G10))>0")
$newG.Priority = 8

$newH = $ws2.Range("H10:L1000000").FormatConditions.Add(2, [System.Reflection.Missing]::Value, "=LEN(TRIM(H10))>0")
$newH.Priority = 7

$newL = $ws2.Range("L9:L1000000").FormatConditions.Add(2, [System.Reflection.Missing]::Value, "=LEN(TRIM(L9))>0")
$newL.Priority = 9
